$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO" (sheet1)
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M32").Value = 712.76
$ws1.Range("M35").Value = 1329.6
$ws1.Range("D56").Value = 1451.52
$ws1.Range("M56").Value = 1568.12
$ws1.Range("D61").Value = "8 de 59"
$ws1.Range("M61").Value = "13 de 59"

# Sheet "VENTA MENSUAL" (sheet2)
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F32").Value = 712.76
$ws2.Range("F35").Value = 1329.6
$ws2.Range("F56").Value = 3019.64
$ws2.Range("F61").Value = 36074.41

# Sheet "CUMPLIMIENTO MENSUAL" (sheet3)
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 11367.15
$ws3.Range("E3").Value = 3458.26
$ws3.Range("F3").Value = 0.7667342758142945

$ws3.Range("D12").Value = 15301.63
$ws3.Range("E12").Value = 35005.37
$ws3.Range("F12").Value = 0.3041650267358419

$ws3.Range("D14").Value = 39040.35
$ws3.Range("E14").Value = 58821.53766749099
$ws3.Range("F14").Value = 0.3989331386356337
